# Add new column 'Servised by' (column O) to the Card24 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

$lastRow = 12

# For every data row, the existing column N cell is currently blank
# (t="inlineStr" with no text). Copy that blank cell's formatting into the
# new column O first (so O ends up as an equivalent "present but empty"
# cell), then fill column N with its new "nan" text.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range("N$r").Copy()
    $ws.Range("O$r").PasteSpecial(-4122)
    $ws.Range("N$r").Value = "nan"
}

# Build the new header cell in O1, copying the header formatting
# (bold, centered, bordered) from the neighbouring N1 header cell.
$ws.Range("N1").Copy()
$ws.Range("O1").PasteSpecial(-4122)
$ws.Range("O1").Value = "Servised by"

$excel.CutCopyMode = 0
